$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.11767664784914444
$ws.Range("A2").Value = -0.0059999999148061534
$ws.Range("A3").Value = -0.003999999934416465
$ws.Range("A4").Value = -0.0079999998774216152
$ws.Range("A5").Value = -0.0029999999465708527
$ws.Range("A6").Value = 0.082239552746687039
$ws.Range("A7").Value = -0.0099999998515083455
$ws.Range("A8").Value = -0.0099999998540463153
$ws.Range("A9").Value = -0.0019999999725892614
$ws.Range("A10").Value = -0.0019999999801179058
$ws.Range("A11").Value = -0.0029999999668746113
$ws.Range("A12").Value = 0.032572844116375954
$ws.Range("A13").Value = -0.0034999999618667843
$ws.Range("A14").Value = -0.0079999998998134814
$ws.Range("A15").Value = -0.00099999999860056477
$ws.Range("A16").Value = 0.011904650977811126
$ws.Range("A17").Value = -0.0019999999855588868
$ws.Range("A18").Value = -0.0015112604061275903
$ws.Range("A19").Value = -0.0039999999419220167
$ws.Range("A20").Value = -0.020282035066445658
$ws.Range("A21").Value = -0.0039999999303939049
$ws.Range("A22").Value = -0.003999999929812148
$ws.Range("A23").Value = -0.0049999999156407071
$ws.Range("A24").Value = -0.019999999698339543
$ws.Range("A25").Value = -0.019999999694007009
$ws.Range("A26").Value = -0.0024999999504675685
$ws.Range("A27").Value = -0.002499999944063358
$ws.Range("A28").Value = -0.001999999923115503
$ws.Range("A29").Value = -0.0069999998345560144
$ws.Range("A30").Value = -0.059999999086601363
$ws.Range("A31").Value = -0.0069999998203122971
$ws.Range("A32").Value = -0.0099999997778450478
$ws.Range("A33").Value = -0.0039999998614757004
